$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: USDT ---
$ws.Range("B2").Value = "trade"
$ws.Range("C2").Value = 1326.50294401
$ws.Range("D2").Value = 795.32972885
$ws.Range("E2").Value = 531.17321516
$ws.Range("G2").Value = 1326.5

# --- Row 3: USDC ---
$ws.Range("B3").Value = "trade"
$ws.Range("C3").Value = 1223.29700881
$ws.Range("D3").Value = 1223.29700881
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 1223.3

# --- Row 4: BTC ---
$ws.Range("B4").Value = "trade"
$ws.Range("C4").Value = 0.02448253
$ws.Range("D4").Value = 0.02448253
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 44344.6
$ws.Range("G4").Value = 1085.67

# --- Row 5: now BTC3S (was ATOM) ---
$ws.Range("A5").Value = "BTC3S"
$ws.Range("B5").Value = "trade"
$ws.Range("C5").Value = 116.1203
$ws.Range("D5").Value = 1.1076
$ws.Range("E5").Value = 115.0127
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 116.12

# --- Row 6: now BTC3L (was ALGO) - new currency ---
$ws.Range("A6").Value = "BTC3L"
$ws.Range("B6").Value = "trade"
$ws.Range("C6").Value = 114.3699
$ws.Range("D6").Value = 3.8651
$ws.Range("E6").Value = 110.5048
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 114.37

# --- Row 7: now ATOM (was BTC3S) ---
$ws.Range("A7").Value = "ATOM"
$ws.Range("B7").Value = "trade"
$ws.Range("C7").Value = 0.2366
$ws.Range("D7").Value = 0.2366
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 13.402
$ws.Range("G7").Value = 3.17

# --- Row 8: now ALGO (was ETH) ---
$ws.Range("A8").Value = "ALGO"
$ws.Range("B8").Value = "trade"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.8378
$ws.Range("G8").Value = 0

# --- Row 9: now ETH (was MATIC) ---
$ws.Range("A9").Value = "ETH"
$ws.Range("B9").Value = "trade"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 3054.33
$ws.Range("G9").Value = 0

# --- Row 10: new row, MATIC ---
$ws.Range("A10").Value = "MATIC"
$ws.Range("B10").Value = "trade"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0
